# Data column cleanup: replace spaces with underscores / remove spaces
# in specific categorical values, per commit message:
# "first sweep cleaning data columns to conform to specs--done by chase"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value()
        if ($val -eq "Environmental Perturbation") {
            $cell.Value = "Environmental_Perturbation"
        } elseif ($val -eq "KN99 alpha") {
            $cell.Value = "KN99_alpha"
        } elseif ($val -eq "Time Course") {
            $cell.Value = "Timecourse"
        }
    }
}

# Update the visible sheet view / selection state to match the authored edit:
# scroll so column E is the left-most visible column (topLeftCell "E1"),
# and select N3 as the active cell (single-cell selection, replacing the
# previous B2:B37 range selection).
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 5
$ws.Range("N3").Select()
